$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3377.8823
$ws.Range("I62").Value = 3461.3
$ws.Range("J62").Value = 3258.7144
$ws.Range("K62").Value = 3461.3
$ws.Range("L62").Value = 3258.7144
$ws.Range("M62").Value = -2837.3
$ws.Range("N62").Value = -4506.7144

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3377.8823
$ws.Range("I65").Value = 3461.3
$ws.Range("J65").Value = 3258.7144
$ws.Range("K65").Value = 17306.5
$ws.Range("L65").Value = 16293.572
$ws.Range("M65").Value = -14186.5
$ws.Range("N65").Value = -22533.572

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 511.26666
$ws.Range("I92").Value = 534.44446
$ws.Range("J92").Value = 476.5
$ws.Range("K92").Value = 534.44446
$ws.Range("L92").Value = 476.5
$ws.Range("M92").Value = 713.55554
$ws.Range("N92").Value = -2972.5

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1596.8
$ws.Range("J112").Value = 1682.2858
$ws.Range("L112").Value = 5046.857400000001
$ws.Range("N112").Value = -7262.857400000001

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3227.3076
$ws.Range("J131").Value = 4215
$ws.Range("L131").Value = 12645
$ws.Range("N131").Value = -22725

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1037463.94
$ws.Range("I138").Value = 2505.5715
$ws.Range("J138").Value = 1361853.9
$ws.Range("K138").Value = 7516.7145
$ws.Range("L138").Value = 4085561.7
$ws.Range("M138").Value = -2376.7145
$ws.Range("N138").Value = -4095841.7

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1973.9062
$ws.Range("I141").Value = 1771.2413
$ws.Range("J141").Value = 3933
$ws.Range("K141").Value = 5313.7239
$ws.Range("L141").Value = 11799
$ws.Range("M141").Value = -133.7239
$ws.Range("N141").Value = -22159

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13690.305
$ws.Range("I32").Value = 10894.76
$ws.Range("K32").Value = 10894.76
$ws.Range("M32").Value = -10607.76

# ARM row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 125043130
$ws.Range("J44").Value = 125043130
$ws.Range("L44").Value = 125043130
$ws.Range("N44").Value = -125044106

# ARM row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 56333.332
$ws.Range("J55").Value = 56333.332
$ws.Range("L55").Value = 56333.332
$ws.Range("N55").Value = -56963.332

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9250.352999999999
$ws.Range("I74").Value = 4480
$ws.Range("J74").Value = 16065.143
$ws.Range("K74").Value = 4480
$ws.Range("L74").Value = 16065.143
$ws.Range("M74").Value = -3606
$ws.Range("N74").Value = -17813.143

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9250.352999999999
$ws.Range("I77").Value = 4480
$ws.Range("J77").Value = 16065.143
$ws.Range("K77").Value = 22400
$ws.Range("L77").Value = 80325.715
$ws.Range("M77").Value = -18032
$ws.Range("N77").Value = -89061.715

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 20836316
$ws.Range("I122").Value = 3298.3333
$ws.Range("K122").Value = 9894.999899999999
$ws.Range("M122").Value = -7444.999899999999

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 5061.4
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 5061.4
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1424.3889
$ws.Range("I107").Value = 1495.5
$ws.Range("J107").Value = 1282.1666
$ws.Range("K107").Value = 1495.5
$ws.Range("L107").Value = 1282.1666
$ws.Range("M107").Value = 424.5
$ws.Range("N107").Value = -5122.1666

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 24307.426
$ws.Range("I134").Value = 2751.6287
$ws.Range("K134").Value = 8254.8861
$ws.Range("M134").Value = -5719.8861

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3763.625
$ws.Range("I31").Value = 1003.5833
$ws.Range("J31").Value = 6523.6665
$ws.Range("K31").Value = 1003.5833
$ws.Range("L31").Value = 6523.6665
$ws.Range("M31").Value = -708.5833
$ws.Range("N31").Value = -7113.6665

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3763.625
$ws.Range("I34").Value = 1003.5833
$ws.Range("J34").Value = 6523.6665
$ws.Range("K34").Value = 1003.5833
$ws.Range("L34").Value = 6523.6665
$ws.Range("M34").Value = -801.5833
$ws.Range("N34").Value = -6927.6665

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2334204.8
$ws.Range("I58").Value = 3790362.5
$ws.Range("K58").Value = 3790362.5
$ws.Range("M58").Value = -3790159.5

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2599.7368
$ws.Range("I134").Value = 1598.5
$ws.Range("J134").Value = 3712.2222
$ws.Range("K134").Value = 4795.5
$ws.Range("L134").Value = 11136.6666
$ws.Range("M134").Value = -2260.5
$ws.Range("N134").Value = -16206.6666

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2334204.8
$ws.Range("I136").Value = 3790362.5
$ws.Range("K136").Value = 11371087.5
$ws.Range("M136").Value = -11368537.5

# CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 104140.1
$ws.Range("I56").Value = 104140.1
$ws.Range("K56").Value = 104140.1
$ws.Range("M56").Value = -103610.1

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1248.6111
$ws.Range("I122").Value = 604.5
$ws.Range("J122").Value = 1570.6666
$ws.Range("K122").Value = 5440.5
$ws.Range("L122").Value = 14135.9994
$ws.Range("M122").Value = -2990.5
$ws.Range("N122").Value = -19035.9994

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 544.08
$ws.Range("I131").Value = 305.67307
$ws.Range("J131").Value = 802.3542
$ws.Range("K131").Value = 917.0192099999999
$ws.Range("L131").Value = 2407.0626
$ws.Range("M131").Value = 4122.98079
$ws.Range("N131").Value = -12487.0626

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1490.9131
$ws.Range("I132").Value = 1265.8889
$ws.Range("K132").Value = 11393.0001
$ws.Range("M132").Value = -8863.000099999999

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4118.5
$ws.Range("I113").Value = 4118.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4118.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1948.5
$ws.Range("N113").ClearContents()

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8407.5
$ws.Range("I136").Value = 8671.429
$ws.Range("K136").Value = 26014.287
$ws.Range("M136").Value = -23464.287

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1754.5454
$ws.Range("I126").Value = 1657.1428
$ws.Range("J126").Value = 1925
$ws.Range("K126").Value = 4971.428400000001
$ws.Range("L126").Value = 5775
$ws.Range("M126").Value = -2501.428400000001
$ws.Range("N126").Value = -10715

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 125044264
$ws.Range("J135").Value = 125044264
$ws.Range("L135").Value = 125044264
$ws.Range("N135").Value = -125054404

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4159.5835
$ws.Range("I136").Value = 2126.827
$ws.Range("J136").Value = 7462.8125
$ws.Range("K136").Value = 6380.481000000001
$ws.Range("L136").Value = 22388.4375
$ws.Range("M136").Value = -3830.481000000001
$ws.Range("N136").Value = -27488.4375
